$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '30.709.49'
Set-TextValue 'E2' '  +1.64%  '

Set-TextValue 'D3' '1.891.99'
Set-TextValue 'E3' '  +1.12%  '

Set-TextValue 'D5' '245.18'
Set-TextValue 'E5' '  +4.51%  '

Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  +0.02%  '

Set-TextValue 'D7' '0.4789'
Set-TextValue 'E7' '  +1.94%  '

Set-TextValue 'E8' '  +1.91%  '

Set-TextValue 'B9' 'Dogecoin'
Set-TextValue 'C9' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D9' '0.06577'
Set-TextValue 'E9' '  +0.16%  '

Set-TextValue 'B10' 'Solana'
Set-TextValue 'C10' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D10' '21.71'
Set-TextValue 'E10' '  +2.15%  '

Set-TextValue 'B11' 'TRON'
Set-TextValue 'C11' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D11' '0.07783'
Set-TextValue 'E11' '  +0.00%  '

Set-TextValue 'B12' 'Litecoin'
Set-TextValue 'C12' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D12' '97.73'
Set-TextValue 'E12' '  +1.51%  '

Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.903.86'
Set-TextValue 'E13' '  +1.76%  '

Set-TextValue 'B14' 'Polygon'
Set-TextValue 'C14' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D14' '0.7451'
Set-TextValue 'E14' '  +8.03%  '

Set-TextValue 'B15' 'Polkadot'
Set-TextValue 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D15' '5.189'
Set-TextValue 'E15' '  +1.93%  '

Set-TextValue 'B16' 'BitcoinCash'
Set-TextValue 'C16' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D16' '281.30'
Set-TextValue 'E16' '  +5.40%  '

Set-TextValue 'B17' 'WrappedBTC'
Set-TextValue 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D17' '30.708.58'
Set-TextValue 'E17' '  +1.67%  '

Set-TextValue 'B18' 'Avalanche'
Set-TextValue 'C18' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D18' '13.51'
Set-TextValue 'E18' '  -1.27%  '

Set-TextValue 'B19' 'ShibaInu'
Set-TextValue 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.000007627'
Set-TextValue 'E19' '  -1.13%  '

Set-TextValue 'B20' 'Dai'
Set-TextValue 'C20' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D20' '1.0000'
Set-TextValue 'E20' '  +0.02%  '

Set-TextValue 'B21' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D21' '2.146.57'
Set-TextValue 'E21' '  +1.37%  '

Set-TextValue 'B22' 'Uniswap'
Set-TextValue 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D22' '5.294'
Set-TextValue 'E22' '  +1.21%  '

Set-TextValue 'B23' 'BinanceUSD'
Set-TextValue 'C23' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D23' '1.001'
Set-TextValue 'E23' '  +0.09%  '

Set-TextValue 'B24' 'Chainlink'
Set-TextValue 'C24' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D24' '6.242'
Set-TextValue 'E24' '  +1.24%  '

Set-TextValue 'B25' 'Cosmos'
Set-TextValue 'C25' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D25' '9.357'
Set-TextValue 'E25' '  -1.19%  '

Set-TextValue 'B26' 'Monero'
Set-TextValue 'C26' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '166.00'
Set-TextValue 'E26' '  +0.30%  '

Set-TextValue 'B27' 'EthereumClassic'
Set-TextValue 'C27' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D27' '19.16'
Set-TextValue 'E27' '  +2.20%  '

Set-TextValue 'B28' 'LidoDAOToken'
Set-TextValue 'C28' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D28' '1.967'
Set-TextValue 'E28' '  +1.61%  '

Set-TextValue 'B29' 'Toncoin'
Set-TextValue 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D29' '1.372'
Set-TextValue 'E29' '  +0.15%  '

Set-TextValue 'B30' 'Stellar'
Set-TextValue 'C30' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D30' '0.09986'
Set-TextValue 'E30' '  +0.58%  '

Set-TextValue 'B31' 'PancakeSwap'
Set-TextValue 'C31' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D31' '1.514'
Set-TextValue 'E31' '  +3.88%  '

Set-TextValue 'B32' 'Filecoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '4.369'
Set-TextValue 'E32' '  +0.27%  '

Set-TextValue 'B33' 'InternetComputer(DFINITY)'
Set-TextValue 'C33' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D33' '4.134'
Set-TextValue 'E33' '  +2.11%  '

Set-TextValue 'B34' 'Hedera'
Set-TextValue 'C34' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D34' '0.04778'
Set-TextValue 'E34' '  +1.01%  '

Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '1.132'
Set-TextValue 'E35' '  +0.25%  '

Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.7045'
Set-TextValue 'E36' '  +0.66%  '

Set-TextValue 'B37' 'HuobiToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D37' '2.719'
Set-TextValue 'E37' '  +0.08%  '

Set-TextValue 'B38' 'VeChain'
Set-TextValue 'C38' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D38' '0.01879'
Set-TextValue 'E38' '  +0.85%  '

Set-TextValue 'B39' 'MXToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D39' '2.765'
Set-TextValue 'E39' '  -0.54%  '

Set-TextValue 'B40' 'FraxShare'
Set-TextValue 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '6.419'
Set-TextValue 'E40' '  +2.93%  '

Set-TextValue 'B41' 'Aave'
Set-TextValue 'C41' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D41' '70.60'
Set-TextValue 'E41' '  -2.86%  '

Set-TextValue 'B42' 'TheSandbox'
Set-TextValue 'C42' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D42' '0.4224'
Set-TextValue 'E42' '  +1.89%  '

Set-TextValue 'B43' 'RenderToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D43' '1.930'
Set-TextValue 'E43' '  -0.39%  '

Set-TextValue 'B44' 'TrustWalletToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D44' '0.8446'
Set-TextValue 'E44' '  +1.31%  '

Set-TextValue 'B45' 'PaxDollar'
Set-TextValue 'C45' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D45' '1.000'
Set-TextValue 'E45' '  +0.00%  '

Set-TextValue 'B46' 'Quant'
Set-TextValue 'C46' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D46' '102.54'
Set-TextValue 'E46' '  -0.22%  '

Set-TextValue 'B47' 'EnergySwap'
Set-TextValue 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '9.411'
Set-TextValue 'E47' '  +2.45%  '

Set-TextValue 'B48' 'Aptos'
Set-TextValue 'C48' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D48' '7.160'
Set-TextValue 'E48' '  +1.35%  '

Set-TextValue 'B49' 'Maker'
Set-TextValue 'C49' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D49' '934.28'
Set-TextValue 'E49' '  -3.85%  '

Set-TextValue 'B50' 'Elrond'
Set-TextValue 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D50' '35.39'
Set-TextValue 'E50' '  +2.65%  '

Set-TextValue 'B51' 'Decentraland'
Set-TextValue 'C51' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D51' '0.3889'
Set-TextValue 'E51' '  +1.86%  '
